{"js": "// Rename nomenclature in \"Justifikasi / bukti ...\" cell:\n// \"Justifikasi / bukti rencana usaha dan/atau kegiatan secara prinsip dapat dilakukan\"\n//   -> \"Justifikasi / bukti persetujuan awal rencana usaha dan/atau kegiatan\"\nconst body = context.document.body;\n\nconst results = body.search(\n  \"Justifikasi / bukti rencana usaha dan/atau kegiatan secara prinsip dapat dilakukan\",\n  { matchCase: true, matchWholeWord: false }\n);\nresults.load(\"text\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\n    \"Justifikasi / bukti persetujuan awal rencana usaha dan/atau kegiatan\",\n    Word.InsertLocation.replace\n  );\n}\n\nawait context.sync();\n", "ps1": "# Rename nomenclature in \"Justifikasi / bukti ...\" cell:\n# \"Justifikasi / bukti rencana usaha dan/atau kegiatan secara prinsip dapat dilakukan\"\n#   -> \"Justifikasi / bukti persetujuan awal rencana usaha dan/atau kegiatan\"\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"Justifikasi / bukti rencana usaha dan/atau kegiatan secara prinsip dapat dilakukan\"\n$find.Replacement.Text = \"Justifikasi / bukti persetujuan awal rencana usaha dan/atau kegiatan\"\n$find.Execute([ref]$find.Text, $true, $true, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n"}
